$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column in H
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
